$d = $word.ActiveDocument

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r") -eq "PI (Tech)") {
        $para.Range.Delete()
        break
    }
}
